$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 4, 6, 7, 9, 10)

foreach ($r in $rows) {
    $ws.Range("B$r").Value = "No"
    $ws.Range("F$r").Value = $null
    $ws.Range("J$r").Value = 0
}
